$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Item numbers for the new/continued rows first
$ws.Range("B33").Value = 10
$ws.Range("B34").Value = 11
$ws.Range("B35").Value = 12
$ws.Range("B36").Value = 13
$ws.Range("B37").Value = 14

# Fill in the description/participant cells in the same order the
# original author typed them (this governs the shared-string table order)
$ws.Range("C33").Value = "Se investiga sobre apis"
$ws.Range("C34").Value = "Se implementa api"
$ws.Range("D35").Value = "Damian Valderrama- Sebastian Espinoza"
$ws.Range("C35").Value = "Se mejora html galeria y api"
$ws.Range("C37").Value = "Se graba el video"
$ws.Range("D34").Value = "Damian Valderrama- Sebastian Espinoza - Nicolas Venegas"
$ws.Range("D37").Value = "Nicolas Venegas- Damian Valderrama - Sebastian Espinoza"
$ws.Range("C36").Value = "Se agrega nuevo implemento en formulario registrarse"

# Remaining cells that reuse existing shared strings
$ws.Range("D33").Value = "Damian Valderrama - Sebastian Espinoza - Nicolas Venegas"
$ws.Range("D36").Value = "Damian Valderrama"

# Apply the same centered style used by the other "item" column cells
# (xlCenter = -4108)
$ws.Range("B34:B37").HorizontalAlignment = -4108

# Update the view to match the target selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("C36").Select()

$wb.Save()
